$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.407.32'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.59%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.579.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.80%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.53%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("E8").Value = '  +0.55%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.80'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.74%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.100'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.53%  '

$ws.Range("E11").Value = '  -1.90%  '

$ws.Range("E12").Value = '  +2.97%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.041.07'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '58.383.90'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.50%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.68'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.77%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.599.51'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.88%  '

$ws.Range("E17").Value = '  -2.12%  '

$ws.Range("E18").Value = '  +0.97%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '335.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.90%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.419'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.83%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.30%  '

$ws.Range("E26").Value = '  -4.86%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0735'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("E30").Value = '  -1.15%  '

$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.96'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.02%  '

$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.94'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.83'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.52%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.89'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.13%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '36.95'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.88%  '

$ws.Range("E36").Value = '  -5.41%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.844'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.818'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.41'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.59'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.47%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '279.75'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.90%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.06%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.590'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.49%  '

$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.64'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.94%  '

$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0533'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.38%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0938'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '18.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.70%  '

$ws.Range("E48").Value = '  -0.49%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.915.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.67%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.40'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.88%  '
